$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 36

# Write the date as plain text (matches existing "MM/DD/YYYY" text cells in column A).
# Forcing a Text number format first prevents Excel from auto-converting the
# string into a date serial number; resetting the Style back to "Normal"
# afterwards avoids leaving a stray cell-style behind.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "10/07/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.151313202125494
$ws.Cells.Item($row, 3).Value = 0.848686797874506
